# Add the MySQL stored-procedure snippet (and its "PROCEDURE" column header)
# to the "MySQL" worksheet, and make that sheet the active tab.
# ("add PROCEDURE :)))) 1400/02/05")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MySQL")

$procedureSql = "DROP PROCEDURE IF EXISTS ``check_login``;`nCREATE PROCEDURE ``check_login``(IN email VARCHAR(200) , IN pass VARCHAR(200))`nBEGIN`n  SELECT * FROM users`n  WHERE users.email COLLATE utf8_unicode_ci = email`n AND`n users.password COLLATE utf8_unicode_ci = pass`n ;`nEND;`nCALL check_login('hassan1@gmail.com','123456')"

$ws.Range("A1").Value = $procedureSql
$ws.Range("B1").Value = "PROCEDURE"

# Formatting: A1 wraps text, B1 is centered both ways (reuses existing styles)
$ws.Range("A1").WrapText = $true
$ws.Range("B1").VerticalAlignment = -4108
$ws.Range("B1").HorizontalAlignment = -4108

# Row / column sizing
$ws.Rows.Item(1).RowHeight = 191.25
$ws.Columns.Item(1).ColumnWidth = 74.6
$ws.Columns.Item(2).ColumnWidth = 10.6

# Make MySQL the active sheet / selection, as in the saved workbook
$ws.Range("M20").Select() | Out-Null
